$d = $word.ActiveDocument

# Fix the spelling typo "oonly" -> "only" in the MAT 337 bullet paragraph.
# This also has the effect of consolidating the surrounding runs/proofErr
# markup when Word re-writes the paragraph text.
$d.Content.Find.Execute("oonly", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "only", 2)
